# CIERRE 11 FEB 2022
# Update the payroll slip ("RECIBO DE NOMINA") for the new pay period and
# the individual figures that changed for JOSE LEOPOLDO ALVARADO GARCIA.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Week banner text (shared string referenced by B9; H9/B27/H27/B43/
#        H43/B60 all pick it up through their existing formulas) ---------
$ws.Range("B9").Value = "SEMANA   06  DEL    7      Al   13   DE   FEBRERO          2022"

# --- 2. JOSE LEOPOLDO ALVARADO GARCIA block (rows 21-26) -----------------
# Dias column: 2 -> 5 (E21 = 500*D21 recalculates automatically)
$ws.Range("D21").Value = 5
# Right-hand "K" total column, plain value (no formula)
$ws.Range("K21").Value = 0
# I S R  (E22)
$ws.Range("E22").Value = -161.65
# PRESTAMO (E25)
$ws.Range("E25").Value = -283
# K24 (=SUM(K21:K23)) and E26 (=SUM(E21:E25)) recalc automatically.

# --- 3. View state: scroll position + active selection -------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
[void]$ws.Range("I13").Select()

# --- 4. Force a full recalculation so every dependent formula (including
#        the TODAY()-based closing-date cells C14/I14/C32/I32/C48/I48/C65)
#        picks up the refreshed cached values. ---------------------------
[void]$excel.CalculateFull()
